$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> take former row 5 values (date 2022-07-07)
$ws.Range("D2").Value = 44749
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17556
$ws.Range("P2").Value = 1170

# Row 4 -> take former row 2 values (date 2023-05-17)
$ws.Range("D4").Value = 45063
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21500
$ws.Range("P4").Value = 1433

# Row 5 -> take former row 4 values (date 2022-10-05)
$ws.Range("D5").Value = 44839
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15600
$ws.Range("P5").Value = 1040

# New row 6 -> new weekly entry (date 2023-06-07)
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 45084
$ws.Range("D6").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112035
$ws.Range("G6").Value = "Bruselas (repollito)"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 22000
$ws.Range("L6").Value = 23000
$ws.Range("M6").Value = 22556
$ws.Range("N6").Value = "$/malla 15 kilos"
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1504
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = "Hortaliza"
